$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the two workers that dropped out of this period's data (delete
# from the bottom up so row numbers above the deleted row stay valid) ---
$ws.Rows(29).Delete()   # 73113114 NICOLAS JOSE GONZALEZ NAVAS
$ws.Rows(20).Delete()   # 73572579 WILFREDO HEREDIA OCHOA

# --- Update the summary figures at the top of the statement ---
$ws.Range("E11").Value = 541284
$ws.Range("C13").Value = 10

# --- Rewrite the worker detail table (rows 16-28) with the refreshed data ---
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "9147896"
$ws.Range("D16").Value = "RICARDO VALDELAMAR LOPEZ"
$ws.Range("E16").Value = "1811"
$ws.Range("F16").Value = 31249
$ws.Range("G16").Value = 781242

$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "12596344"
$ws.Range("D17").Value = "DOMINGO SANTANA TORRES ALEMAN"
$ws.Range("E17").Value = "1811"
$ws.Range("F17").Value = 31249
$ws.Range("G17").Value = 781242

$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "92400728"
$ws.Range("D18").Value = "PEDRO LUIS GALAN MERCADO"
$ws.Range("E18").Value = "1811"
$ws.Range("F18").Value = 31249
$ws.Range("G18").Value = 781242

$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "73581352"
$ws.Range("D19").Value = "FELIPE MONTERROSA CAMARGO"
$ws.Range("E19").Value = "1812"
$ws.Range("F19").Value = 31249
$ws.Range("G19").Value = 781242

$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "9147896"
$ws.Range("D20").Value = "RICARDO VALDELAMAR LOPEZ"
$ws.Range("E20").Value = "1812"
$ws.Range("F20").Value = 31249
$ws.Range("G20").Value = 781242

$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "1143348532"
$ws.Range("D21").Value = "LUIS FELIPE LEON JIMENEZ"
$ws.Range("E21").Value = "1812"
$ws.Range("F21").Value = 67200
$ws.Range("G21").Value = 1680000

$ws.Range("B22").Value = "CC"
$ws.Range("C22").Value = "1047397104"
$ws.Range("D22").Value = "YASANDRA MILENA SANTOS PEREZ"
$ws.Range("E22").Value = "1812"
$ws.Range("F22").Value = 37490
$ws.Range("G22").Value = 937250

$ws.Range("B23").Value = "CC"
$ws.Range("C23").Value = "73084660"
$ws.Range("D23").Value = "RAUL TADEO LEON LOPEZ"
$ws.Range("E23").Value = "1812"
$ws.Range("F23").Value = 84000
$ws.Range("G23").Value = 2100000

$ws.Range("B24").Value = "CC"
$ws.Range("C24").Value = "12596344"
$ws.Range("D24").Value = "DOMINGO SANTANA TORRES ALEMAN"
$ws.Range("E24").Value = "1812"
$ws.Range("F24").Value = 31249
$ws.Range("G24").Value = 781242

$ws.Range("B25").Value = "CC"
$ws.Range("C25").Value = "8852958"
$ws.Range("D25").Value = "GUSTAVO BALLESTAS CASTILLO"
$ws.Range("E25").Value = "1812"
$ws.Range("F25").Value = 31249
$ws.Range("G25").Value = 781242

$ws.Range("B26").Value = "CC"
$ws.Range("C26").Value = "1047375080"
$ws.Range("D26").Value = "SAMIR ENRIQUE JULIAO DAGER"
$ws.Range("E26").Value = "1812"
$ws.Range("F26").Value = 70602
$ws.Range("G26").Value = 1765050

$ws.Range("B27").Value = "CC"
$ws.Range("C27").Value = "92400728"
$ws.Range("D27").Value = "PEDRO LUIS GALAN MERCADO"
$ws.Range("E27").Value = "1812"
$ws.Range("F27").Value = 31249
$ws.Range("G27").Value = 781242

$ws.Range("B28").Value = "CC"
$ws.Range("C28").Value = "1118859754"
$ws.Range("D28").Value = "JESUS JOSE IBARRA TORRENEGRA"
$ws.Range("E28").Value = "1812"
$ws.Range("F28").Value = 32000
$ws.Range("G28").Value = 800000
